$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update summary figures (row 11 "VALOR MORA", row 13 counts) ---
$ws.Range("E11").Value = 268693
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 5

# --- Give the about-to-become-last data row (21) the "closing" border
#     formatting that currently belongs to row 22 (the last row of the
#     table, which is being removed). Copy formats only, before the
#     content gets overwritten below. ---
$ws.Range("B22:J22").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)

# --- Rewrite the worker detail table (rows 16-21), dropping the
#     "LUIS ROBERTO RAMIREZ NAVARRO" record and adding the new
#     "JORGE LUIS MARTINEZ ORTIZ" / "CARLOS EDUARDO VEGAS MACIAS" periods ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73581883"
$ws.Range("D16").Value = "JORGE LUIS MARTINEZ ORTIZ"
$ws.Range("E16").Value = "2209"
$ws.Range("F16").Value = 38666
$ws.Range("G16").Value = 1000000

$ws.Range("B17").Value = "PPT"
$ws.Range("C17").Value = "4865522"
$ws.Range("D17").Value = "CARLOS EDUARDO VEGAS MACIAS"
$ws.Range("E17").Value = "2301"
$ws.Range("F17").Value = 60000
$ws.Range("G17").Value = 1747200

$ws.Range("B18").Value = "PPT"
$ws.Range("C18").Value = "4865522"
$ws.Range("D18").Value = "CARLOS EDUARDO VEGAS MACIAS"
$ws.Range("E18").Value = "2302"
$ws.Range("F18").Value = 60000
$ws.Range("G18").Value = 1747200

$ws.Range("B19").Value = "PPT"
$ws.Range("C19").Value = "4865522"
$ws.Range("D19").Value = "CARLOS EDUARDO VEGAS MACIAS"
$ws.Range("E19").Value = "2303"
$ws.Range("F19").Value = 60000
$ws.Range("G19").Value = 1747200

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "33104272"
$ws.Range("D20").Value = "NISIDA ALCAZAR JIMENEZ"
$ws.Range("E20").Value = "2309"
$ws.Range("F20").Value = 3627
$ws.Range("G20").Value = 1360000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1043962314"
$ws.Range("D21").Value = "BANIC HAISAWA OTERO LABRADOR"
$ws.Range("E21").Value = "2309"
$ws.Range("F21").Value = 46400
$ws.Range("G21").Value = 1160000

# --- Row 22 (old "LUIS ROBERTO RAMIREZ NAVARRO" line) is removed entirely.
#     This also shifts the trailing signature block up: the old row 27
#     ("___") becomes row 26, and the old row 28 (names/legal text)
#     becomes row 27 - exactly the desired final layout. ---
$ws.Rows.Item(22).Delete()
